# dataDriven.xlsx - 10 Mar 2025 update
#
# The login-credentials table on "LoginSheet" used to list two accounts
# (row 2: pankaj_kalra@unifyedqa.edu / Admin@2008s, row 3: philip_parker@
# unifyedqa.edu / Admin@2008). The pankaj_kalra row was removed, so the
# philip_parker row becomes the new (and only) data row, row 2. The
# saved cell selection also moved from H11 to D11.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the old row 2 (pankaj_kalra@unifyedqa.edu / Admin@2008s).
# This shifts the old row 3 (philip_parker@unifyedqa.edu / Admin@2008) up
# into row 2 and drops the now-unused shared strings automatically.
$ws.Rows(2).Delete()

# Row deletion does not renumber the worksheet's saved hyperlink list, so
# rebuild it to point at the surviving row.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:philip_parker@unifyedqa.edu") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:Admin@2008") | Out-Null

# Adding a hyperlink restyles the cell; put the original "Hyperlink" cell
# style back on the surviving row.
$ws.Range("A2:B2").Style = "Hyperlink"

# Match the saved cell selection recorded in the workbook.
$ws.Range("D11").Select() | Out-Null
